$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "45546839"
$ws.Range("D16").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E16").Value = "2405"
$ws.Range("F16").Value = 58667
$ws.Range("G16").Value = 2000000

$ws.Range("C17").Value = "45546839"
$ws.Range("D17").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E17").Value = "2404"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

$ws.Range("C18").Value = "45546839"
$ws.Range("D18").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E18").Value = "2403"
$ws.Range("F18").Value = 80000
$ws.Range("G18").Value = 2000000

$ws.Range("C19").Value = "45546839"
$ws.Range("D19").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E19").Value = "2402"
$ws.Range("F19").Value = 80000
$ws.Range("G19").Value = 2000000

$ws.Range("C20").Value = "45546839"
$ws.Range("D20").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E20").Value = "2401"
$ws.Range("F20").Value = 80000
$ws.Range("G20").Value = 2000000

$ws.Range("C21").Value = "45546839"
$ws.Range("D21").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E21").Value = "2312"
$ws.Range("F21").Value = 80000
$ws.Range("G21").Value = 2000000

$ws.Range("C22").Value = "45546839"
$ws.Range("D22").Value = "CLAUDIA PATRICIA MEJIA RAMIREZ"
$ws.Range("E22").Value = "2311"
$ws.Range("F22").Value = 80000
$ws.Range("G22").Value = 2000000

$ws.Range("C23").Value = "1047422454"
$ws.Range("D23").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E23").Value = "2405"
$ws.Range("F23").Value = 58667
$ws.Range("G23").Value = 2000000

$ws.Range("C24").Value = "1047422454"
$ws.Range("D24").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E24").Value = "2404"
$ws.Range("F24").Value = 80000
$ws.Range("G24").Value = 2000000

$ws.Range("C25").Value = "1047422454"
$ws.Range("D25").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E25").Value = "2403"
$ws.Range("F25").Value = 80000
$ws.Range("G25").Value = 2000000

$ws.Range("C26").Value = "1047422454"
$ws.Range("D26").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E26").Value = "2402"
$ws.Range("F26").Value = 80000
$ws.Range("G26").Value = 2000000

$ws.Range("C27").Value = "1047422454"
$ws.Range("D27").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E27").Value = "2401"
$ws.Range("F27").Value = 80000
$ws.Range("G27").Value = 2000000

$ws.Range("C28").Value = "1047422454"
$ws.Range("D28").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E28").Value = "2312"
$ws.Range("F28").Value = 80000
$ws.Range("G28").Value = 2000000

$ws.Range("C29").Value = "1047422454"
$ws.Range("D29").Value = "MARTICELA BAYTER DORIA"
$ws.Range("E29").Value = "2311"
$ws.Range("F29").Value = 80000
$ws.Range("G29").Value = 2000000

$ws.Range("C30").Value = "1143360875"
$ws.Range("D30").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E30").Value = "2405"
$ws.Range("F30").Value = 58667
$ws.Range("G30").Value = 1800000

$ws.Range("C31").Value = "1143360875"
$ws.Range("D31").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E31").Value = "2404"
$ws.Range("F31").Value = 80000
$ws.Range("G31").Value = 1800000

$ws.Range("C32").Value = "1143360875"
$ws.Range("D32").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E32").Value = "2403"
$ws.Range("F32").Value = 80000
$ws.Range("G32").Value = 1800000

$ws.Range("C33").Value = "1143360875"
$ws.Range("D33").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E33").Value = "2402"
$ws.Range("F33").Value = 80000
$ws.Range("G33").Value = 1800000

$ws.Range("C34").Value = "1143360875"
$ws.Range("D34").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E34").Value = "2401"
$ws.Range("F34").Value = 80000
$ws.Range("G34").Value = 1800000

$ws.Range("C35").Value = "1143360875"
$ws.Range("D35").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E35").Value = "2312"
$ws.Range("F35").Value = 80000
$ws.Range("G35").Value = 1800000

$ws.Range("C36").Value = "1143360875"
$ws.Range("D36").Value = "JOSE GUILLERMO ANGULO VIAÑA"
$ws.Range("E36").Value = "2311"
$ws.Range("F36").Value = 80000
$ws.Range("G36").Value = 1800000
